$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.209.68"
$ws.Range("E2").Value = "  -4.14%  "
$ws.Range("D3").Value = "1.956.88"
$ws.Range("E3").Value = "  -4.59%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "'241.88"
$ws.Range("E5").Value = "  -4.04%  "
$ws.Range("E6").Value = "  -5.24%  "
$ws.Range("B7").Value = "Solana"
$ws.Range("C7").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D7").Value = "'57.54"
$ws.Range("E7").Value = "  -11.72%  "
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("D9").Value = "'0.368"
$ws.Range("E9").Value = "  -2.84%  "
$ws.Range("D10").Value = "'56.73"
$ws.Range("E10").Value = "  -6.93%  "
$ws.Range("D11").Value = "'0.0783"
$ws.Range("E11").Value = "  +2.60%  "
$ws.Range("D12").Value = "'0.103"
$ws.Range("E12").Value = "  -1.07%  "
$ws.Range("D13").Value = "'0.838"
$ws.Range("E13").Value = "  -8.77%  "
$ws.Range("D14").Value = "'13.83"
$ws.Range("E14").Value = "  -8.73%  "
$ws.Range("D15").Value = "'21.54"
$ws.Range("E15").Value = "  +4.00%  "
$ws.Range("D16").Value = "2.245.59"
$ws.Range("E16").Value = "  -4.43%  "
$ws.Range("E17").Value = "  -3.60%  "
$ws.Range("D18").Value = "1.955.77"
$ws.Range("E18").Value = "  -4.55%  "
$ws.Range("D19").Value = "36.013.59"
$ws.Range("E19").Value = "  -4.21%  "
$ws.Range("D20").Value = "'70.74"
$ws.Range("E20").Value = "  -4.64%  "
$ws.Range("D21").Value = "0.0₃0844"
$ws.Range("E21").Value = "  -3.78%  "
$ws.Range("D22").Value = "'235.81"
$ws.Range("E22").Value = "  -1.29%  "
$ws.Range("E23").Value = "  -3.74%  "
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("D25").Value = "'2.51"
$ws.Range("E25").Value = "  -6.61%  "
$ws.Range("E26").Value = "  -5.28%  "
$ws.Range("D27").Value = "'9.58"
$ws.Range("E27").Value = "  -0.42%  "
$ws.Range("D28").Value = "'160.32"
$ws.Range("E28").Value = "  +0.20%  "
$ws.Range("D29").Value = "'19.66"
$ws.Range("E29").Value = "  -1.69%  "
$ws.Range("E30").Value = "  +6.00%  "
$ws.Range("E31").Value = "  -2.83%  "
$ws.Range("D32").Value = "'4.82"
$ws.Range("E32").Value = "  -7.74%  "
$ws.Range("D33").Value = "'1.12"
$ws.Range("E33").Value = "  -7.41%  "
$ws.Range("D34").Value = "'0.0607"
$ws.Range("E34").Value = "  -1.91%  "
$ws.Range("E35").Value = "  -7.68%  "
$ws.Range("E36").Value = "  +0.16%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "'2.25"
$ws.Range("E37").Value = "  -7.61%  "
$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").Value = "'1.81"
$ws.Range("E38").Value = "  -2.68%  "
$ws.Range("B39").Value = "THORChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D39").Value = "'5.94"
$ws.Range("E39").Value = "  -2.25%  "
$ws.Range("D40").Value = "'2.96"
$ws.Range("E40").Value = "  +4.71%  "
$ws.Range("D41").Value = "'0.0979"
$ws.Range("E41").Value = "  -5.70%  "
$ws.Range("E42").Value = "  -2.03%  "
$ws.Range("E43").Value = "  -1.75%  "
$ws.Range("E44").Value = "  -3.94%  "
$ws.Range("E45").Value = "  -5.72%  "
$ws.Range("D46").Value = "'90.76"
$ws.Range("E46").Value = "  -4.70%  "
$ws.Range("D47").Value = "'15.74"
$ws.Range("E47").Value = "  -7.39%  "
$ws.Range("D48").Value = "'7.46"
$ws.Range("E48").Value = "  -6.84%  "
$ws.Range("D49").Value = "1.331.81"
$ws.Range("E49").Value = "  -5.64%  "
$ws.Range("E50").Value = "  -4.23%  "
$ws.Range("D51").Value = "2.135.28"
$ws.Range("E51").Value = "  -4.50%  "
